$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data cells for rows 8-11 (columns C, D, E)
$ws.Range("D8").Value = 36848
$ws.Range("E8").Value = 147307

$ws.Range("C9").Value = 44644
$ws.Range("D9").Value = 36854
$ws.Range("E9").Value = 147340

$ws.Range("C10").Value = 44616
$ws.Range("D10").Value = 36829
$ws.Range("E10").Value = 147434

$ws.Range("C11").Value = 44758

# Force recalculation so the AVERAGE formulas in D12/E12 refresh
$excel.Calculate()

# Update the selected cell / range to match the new cursor position
$ws.Range("C11").Select()
